# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.432.05"
Set-TextValue "D3" "1.869.87"
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.14%  "
Set-TextValue "D5" "243.36"
Set-TextValue "E5" "  +0.28%  "
Set-TextValue "D6" "0.7054"
Set-TextValue "E6" "  -2.11%  "
Set-TextValue "D7" "1.000"
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "0.07940"
Set-TextValue "E8" "  -1.02%  "
Set-TextValue "E9" "  -0.25%  "
Set-TextValue "E10" "  -2.13%  "
Set-TextValue "D11" "0.07813"
Set-TextValue "E11" "  -4.58%  "
Set-TextValue "D12" "1.898.66"
Set-TextValue "E12" "  +1.23%  "
Set-TextValue "D13" "93.79"
Set-TextValue "E13" "  -1.09%  "
Set-TextValue "D14" "5.162"
Set-TextValue "E14" "  -1.38%  "
Set-TextValue "D15" "0.7036"
Set-TextValue "E15" "  -1.43%  "
Set-TextValue "D16" "6.488"
Set-TextValue "E16" "  +1.09%  "
Set-TextValue "D17" "0.000008660"
Set-TextValue "E17" "  +1.79%  "
Set-TextValue "D18" "29.519.88"
Set-TextValue "E18" "  +0.55%  "
Set-TextValue "D19" "252.59"
Set-TextValue "E19" "  +3.33%  "
Set-TextValue "D20" "2.159.07"
Set-TextValue "E20" "  +1.68%  "
Set-TextValue "D21" "13.11"
Set-TextValue "E21" "  -1.54%  "
Set-TextValue "E22" "  +0.00%  "
Set-TextValue "D23" "7.665"
Set-TextValue "E23" "  -1.04%  "
Set-TextValue "D24" "1.000"
Set-TextValue "E24" "  -0.30%  "
Set-TextValue "D25" "0.1546"
Set-TextValue "E25" "  -3.05%  "
Set-TextValue "D26" "8.998"
Set-TextValue "E26" "  -0.58%  "
Set-TextValue "D27" "161.42"
Set-TextValue "E27" "  -0.68%  "
Set-TextValue "D28" "18.81"
Set-TextValue "E28" "  +1.55%  "
Set-TextValue "D29" "1.500"
Set-TextValue "E29" "  -0.30%  "
Set-TextValue "D30" "4.307"
Set-TextValue "E30" "  -2.36%  "
Set-TextValue "D31" "4.267"
Set-TextValue "E31" "  -0.50%  "
Set-TextValue "D32" "1.211"
Set-TextValue "E32" "  -1.44%  "
Set-TextValue "D33" "0.05291"
Set-TextValue "E33" "  -1.29%  "
Set-TextValue "D34" "1.907"
Set-TextValue "E34" "  -1.74%  "
Set-TextValue "D35" "0.7615"
Set-TextValue "E35" "  -0.80%  "
Set-TextValue "E36" "  +0.43%  "
Set-TextValue "D37" "2.705"
Set-TextValue "E37" "  +0.09%  "
Set-TextValue "D38" "0.01879"
Set-TextValue "E38" "  +0.31%  "
Set-TextValue "D39" "1.282.81"
Set-TextValue "E39" "  +0.55%  "
Set-TextValue "D40" "2.765"
Set-TextValue "E40" "  +0.44%  "
Set-TextValue "D41" "0.9015"
Set-TextValue "E41" "  -0.85%  "
Set-TextValue "D42" "109.77"
Set-TextValue "E42" "  -3.02%  "
Set-TextValue "D43" "6.016"
Set-TextValue "E43" "  -6.57%  "
Set-TextValue "D44" "70.79"
Set-TextValue "E44" "  -4.84%  "
Set-TextValue "E45" "  -0.14%  "
Set-TextValue "D46" "2.054.95"
Set-TextValue "E46" "  +1.73%  "
Set-TextValue "D47" "0.00000000126"
Set-TextValue "E47" "  -3.60%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.630"
Set-TextValue "E48" "  +1.31%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.804"
Set-TextValue "E49" "  -0.04%  "
Set-TextValue "D50" "0.5174"
Set-TextValue "E50" "  -1.11%  "
Set-TextValue "D51" "0.4297"
Set-TextValue "E51" "  -1.14%  "
